# Prepend a new weekly price record for "Haba" (Feria Lagunitas de Puerto Montt)
# by inserting a new row at the top of the data block (row 80) and shifting
# every existing record down by one row. The former last record (old row 128)
# ends up at row 129, and the brand-new record is written into row 80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 80; Excel shifts rows 80-128 down
# to 81-129 and carries the existing formatting (e.g. the date style on
# column D) onto the newly inserted row.
$ws.Rows(80).Insert()

# Populate the newly inserted row 80 with the new record's data.
$ws.Cells.Item(80, 1).Value = 4
$ws.Cells.Item(80, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(80, 3).Value = "Los Lagos"
$ws.Cells.Item(80, 4).Value = 45072
$ws.Cells.Item(80, 5).Value = 10
$ws.Cells.Item(80, 6).Value = 100112026
$ws.Cells.Item(80, 7).Value = "Haba"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 80
$ws.Cells.Item(80, 11).Value = 22000
$ws.Cells.Item(80, 12).Value = 22000
$ws.Cells.Item(80, 13).Value = 22000
$ws.Cells.Item(80, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(80, 16).Value = 880
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
